$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '53.362.92'
$ws.Range('E2').Value = '  +3.54%  '
$ws.Range('D3').Value = '3.152.61'
$ws.Range('E3').Value = '  +3.11%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '398.06'
$ws.Range('E5').Value = '  +3.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.80'
$ws.Range('E6').Value = '  +5.34%  '
$ws.Range('E7').Value = '  +0.89%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.608'
$ws.Range('E9').Value = '  +3.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.98'
$ws.Range('E10').Value = '  +5.76%  '
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0871'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('D13').Value = '3.649.58'
$ws.Range('E13').Value = '  +3.11%  '
$ws.Range('E14').Value = '  +2.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.99'
$ws.Range('E15').Value = '  +2.74%  '
$ws.Range('E16').Value = '  +9.02%  '
$ws.Range('D17').Value = '3.154.14'
$ws.Range('E17').Value = '  +3.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.58'
$ws.Range('E18').Value = '  -1.80%  '
$ws.Range('D19').Value = '53.346.55'
$ws.Range('E19').Value = '  +3.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.27'
$ws.Range('E20').Value = '  +3.41%  '
$ws.Range('E21').Value = '  +3.27%  '
$ws.Range('D22').Value = '0.0₃0973'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.76'
$ws.Range('E23').Value = '  +0.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.82'
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.24'
$ws.Range('E25').Value = '  +2.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.16'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.70'
$ws.Range('E27').Value = '  +2.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.35'
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.110'
$ws.Range('E31').Value = '  +2.41%  '
$ws.Range('E32').Value = '  +7.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '37.26'
$ws.Range('E33').Value = '  +7.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0498'
$ws.Range('E34').Value = '  +11.41%  '
$ws.Range('E35').Value = '  +0.46%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.65'
$ws.Range('E37').Value = '  +9.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.80'
$ws.Range('E39').Value = '  +9.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.16'
$ws.Range('E40').Value = '  +11.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.292'
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.41'
$ws.Range('E42').Value = '  +2.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.90'
$ws.Range('E43').Value = '  +1.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '130.89'
$ws.Range('E44').Value = '  +4.57%  '
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.32'
$ws.Range('E46').Value = '  +1.86%  '
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('D48').Value = '2.091.03'
$ws.Range('E48').Value = '  +2.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.38'
$ws.Range('E49').Value = '  -1.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0504'
$ws.Range('E50').Value = '  +21.59%  '
$ws.Range('E51').Value = '  +5.67%  '
